$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the log table with a new entry (I2C addresses added to the control-circuit docs)
$lastRow = 24
$newRow = $lastRow + 1

# Carry the date/number formatting down from the row above (keeps the same cell style,
# instead of inventing a new custom number format)
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A" + $newRow).Value = Get-Date -Year 2017 -Month 1 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("B" + $newRow).Value = "Kontrollschaltung2.pdf"
$ws.Range("C" + $newRow).Value = "Addressen hinzugefügt"

$ws.Range("C" + $newRow).Select()
